$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 107-108; this shifts the former rows 107-119
# down to become rows 109-121 (matching the diff's row-shift pattern).
$ws.Rows("107:108").Insert()

# New row 107 - a weekly "Especial" quality record
$ws.Range("A107").Value = 3
$ws.Range("B107").Value = "Femacal de La Calera"
$ws.Range("C107").Value = "Coquimbo"
$ws.Range("D107").Value = 44476
$ws.Range("E107").Value = 5
$ws.Range("F107").Value = "Fruta"
$ws.Range("G107").Value = 100101
$ws.Range("H107").Value = "Berries"
$ws.Range("I107").Value = 100112025
$ws.Range("J107").Value = "Frutilla"
$ws.Range("K107").Value = "Sin especificar"
$ws.Range("L107").Value = "Especial"
$ws.Range("M107").Value = 58
$ws.Range("N107").Value = 14000
$ws.Range("O107").Value = 14000
$ws.Range("P107").Value = 14000
$ws.Range("Q107").Value = "`$/bandeja 7 kilos"
$ws.Range("R107").Value = "Provincia de Melipilla"
$ws.Range("S107").Value = 2000
$ws.Range("T107").Value = 7

# New row 108 - a weekly "Segunda" quality record
$ws.Range("A108").Value = 3
$ws.Range("B108").Value = "Femacal de La Calera"
$ws.Range("C108").Value = "Coquimbo"
$ws.Range("D108").Value = 44476
$ws.Range("E108").Value = 5
$ws.Range("F108").Value = "Fruta"
$ws.Range("G108").Value = 100101
$ws.Range("H108").Value = "Berries"
$ws.Range("I108").Value = 100112025
$ws.Range("J108").Value = "Frutilla"
$ws.Range("K108").Value = "Sin especificar"
$ws.Range("L108").Value = "Segunda"
$ws.Range("M108").Value = 45
$ws.Range("N108").Value = 10000
$ws.Range("O108").Value = 10000
$ws.Range("P108").Value = 10000
$ws.Range("Q108").Value = "`$/bandeja 7 kilos"
$ws.Range("R108").Value = "Provincia de Melipilla"
$ws.Range("S108").Value = 1429
$ws.Range("T108").Value = 7
